# Update Efnb1-Epha4 LR-pair metrics with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 9.546140333333334
$ws.Range("H2").Value = 28.638421
$ws.Range("I2").Value = 0.587227294878132
$ws.Range("J2").Value = 0.587227294878132
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 97.66285148379045
$ws.Range("R2").Value = 878.965663354114
$ws.Range("S2").Value = 0.2851361657720983
$ws.Range("T2").Value = 0.2851361657720983
# Row 3
$ws.Range("G3").Value = 9.546140333333334
$ws.Range("H3").Value = 28.638421
$ws.Range("I3").Value = 0.587227294878132
$ws.Range("J3").Value = 0.587227294878132
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 83.76916018914876
$ws.Range("R3").Value = 753.922441702339
$ws.Range("S3").Value = 0.2445721866952348
$ws.Range("T3").Value = 0.2445721866952348
# Row 4
$ws.Range("G4").Value = 9.546140333333334
$ws.Range("H4").Value = 28.638421
$ws.Range("I4").Value = 0.587227294878132
$ws.Range("J4").Value = 0.587227294878132
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 19.42402813553067
$ws.Range("R4").Value = 174.816253219776
$ws.Range("S4").Value = 0.05671033378882883
$ws.Range("T4").Value = 0.05671033378882884
# Row 5
$ws.Range("G5").Value = 9.546140333333334
$ws.Range("H5").Value = 28.638421
$ws.Range("I5").Value = 0.587227294878132
$ws.Range("J5").Value = 0.587227294878132
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.2769589874442223
$ws.Range("R5").Value = 2.492630886998
$ws.Range("S5").Value = 0.0008086086219699969
$ws.Range("T5").Value = 0.0008086086219699969
# Row 6
$ws.Range("I6").Value = 0.2496684258894083
$ws.Range("J6").Value = 0.2496684258894083
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 41.52281511861489
$ws.Range("R6").Value = 373.705336067534
$ws.Range("S6").Value = 0.1212298854181075
$ws.Range("T6").Value = 0.1212298854181075
# Row 7
$ws.Range("I7").Value = 0.2496684258894083
$ws.Range("J7").Value = 0.2496684258894083
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("S7").Value = 0.1039835058777403
$ws.Range("T7").Value = 0.1039835058777403
# Row 8
$ws.Range("I8").Value = 0.2496684258894083
$ws.Range("J8").Value = 0.2496684258894083
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 8.258414708117332
$ws.Range("R8").Value = 74.325732373056
$ws.Range("S8").Value = 0.02411124260097311
$ws.Range("T8").Value = 0.02411124260097312
# Row 9
$ws.Range("I9").Value = 0.2496684258894083
$ws.Range("J9").Value = 0.2496684258894083
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.1177532363264445
$ws.Range("R9").Value = 1.059779126938
$ws.Range("S9").Value = 0.0003437919925873847
$ws.Range("T9").Value = 0.0003437919925873847
# Row 10
$ws.Range("G10").Value = 2.210442
$ws.Range("H10").Value = 6.631326
$ws.Range("I10").Value = 0.1359745227725727
$ws.Range("J10").Value = 0.1359745227725727
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 22.614172976876
$ws.Range("R10").Value = 203.527556791884
$ws.Range("S10").Value = 0.06602427101776406
$ws.Range("T10").Value = 0.06602427101776406
# Row 11
$ws.Range("G11").Value = 2.210442
$ws.Range("H11").Value = 6.631326
$ws.Range("I11").Value = 0.1359745227725727
$ws.Range("J11").Value = 0.1359745227725727
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 19.397040429026
$ws.Range("R11").Value = 174.573363861234
$ws.Range("S11").Value = 0.05663154056255281
$ws.Range("T11").Value = 0.05663154056255282
# Row 12
$ws.Range("G12").Value = 2.210442
$ws.Range("H12").Value = 6.631326
$ws.Range("I12").Value = 0.1359745227725727
$ws.Range("J12").Value = 0.1359745227725727
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 4.497701280384
$ws.Range("R12").Value = 40.479311523456
$ws.Range("S12").Value = 0.01313147505313017
$ws.Range("T12").Value = 0.01313147505313017
# Row 13
$ws.Range("G13").Value = 2.210442
$ws.Range("H13").Value = 6.631326
$ws.Range("I13").Value = 0.1359745227725727
$ws.Range("J13").Value = 0.1359745227725727
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.06413081693200001
$ws.Range("R13").Value = 0.5771773523880001
$ws.Range("S13").Value = 0.0001872361391256107
$ws.Range("T13").Value = 0.0001872361391256107
# Row 14
$ws.Range("G14").Value = 0.4410293333333333
$ws.Range("H14").Value = 1.323088
$ws.Range("I14").Value = 0.02712975645988715
$ws.Range("J14").Value = 0.02712975645988715
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 10.23061133333333
$ws.Range("N14").Value = 30.691834
$ws.Range("O14").Value = 0.4855635428718841
$ws.Range("P14").Value = 0.4855635428718841
$ws.Range("Q14").Value = 4.511999695932444
$ws.Range("R14").Value = 40.607997263392
$ws.Range("S14").Value = 0.01317322066391419
$ws.Range("T14").Value = 0.01317322066391419
# Row 15
$ws.Range("G15").Value = 0.4410293333333333
$ws.Range("H15").Value = 1.323088
$ws.Range("I15").Value = 0.02712975645988715
$ws.Range("J15").Value = 0.02712975645988715
$ws.Range("O15").Value = 0.4164864079521221
$ws.Range("P15").Value = 0.4164864079521222
$ws.Range("Q15").Value = 3.870114578465777
$ws.Range("R15").Value = 34.831031206192
$ws.Range("S15").Value = 0.01129917481659428
$ws.Range("T15").Value = 0.01129917481659428
# Row 16
$ws.Range("G16").Value = 0.4410293333333333
$ws.Range("H16").Value = 1.323088
$ws.Range("I16").Value = 0.02712975645988715
$ws.Range("J16").Value = 0.02712975645988715
$ws.Range("M16").Value = 2.034752
$ws.Range("N16").Value = 6.104255999999999
$ws.Range("O16").Value = 0.09657305490303886
$ws.Range("P16").Value = 0.09657305490303887
$ws.Range("Q16").Value = 0.8973853180586665
$ws.Range("R16").Value = 8.076467862528
$ws.Range("S16").Value = 0.002620003460106754
$ws.Range("T16").Value = 0.002620003460106755
# Row 17
$ws.Range("G17").Value = 0.4410293333333333
$ws.Range("H17").Value = 1.323088
$ws.Range("I17").Value = 0.02712975645988715
$ws.Range("J17").Value = 0.02712975645988715
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.02901266666666667
$ws.Range("N17").Value = 0.087038
$ws.Range("O17").Value = 0.001376994272954919
$ws.Range("P17").Value = 0.001376994272954919
$ws.Range("Q17").Value = 0.01279543703822222
$ws.Range("R17").Value = 0.115158933344
$ws.Range("S17").Value = 0.000037357519271926308782744153
$ws.Range("T17").Value = 0.000037357519271926308782744153
